# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45175 (2023-09-06) to 45177 (2023-09-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 189; $row++) {
    $ws.Cells.Item($row, 3).Value = 45177
}
